{"js": "// Bold several text fragments inside the \"Worked as a Frontend programmer in\n// Ateneo's ACMS ... under SAMAHAN System Development org.\" bullet point.\n//\n// 1) \"Ateneo's\"                                    -> bold\n// 2) \" ACMS (Automated Crowd Management System)\"    -> bold (\" for IT Week\" stays regular)\n// 3) \"SAMAHAN System Development org\"                -> bold (\" under \" stays regular)\n// 4) the \".\" right after \"org\"                       -> bold (trailing space stays regular)\n\nconst body = context.document.body;\n\n// 1) \"Ateneo's\" (curly apostrophe, same as stored in the document)\nconst ateneoResults = body.search(\"Ateneo\\u2019s\", { matchCase: true });\nateneoResults.load(\"items\");\nawait context.sync();\nfor (const rng of ateneoResults.items) {\n  rng.font.bold = true;\n}\n\n// 2) \" ACMS (Automated Crowd Management System)\" (leave \" for IT Week\" regular)\nconst acmsResults = body.search(\" ACMS (Automated Crowd Management System)\", { matchCase: true });\nacmsResults.load(\"items\");\nawait context.sync();\nfor (const rng of acmsResults.items) {\n  rng.font.bold = true;\n}\n\n// 3) \"SAMAHAN System Development org\" (leave \" under \" regular)\nconst samahanResults = body.search(\"SAMAHAN System Development org\", { matchCase: true });\nsamahanResults.load(\"items,start\");\nawait context.sync();\nconst samahanRange = samahanResults.items[0];\nsamahanRange.font.bold = true;\n\n// 4) the \".\" that immediately follows \"org\" (leave the trailing space regular)\nconst afterSamahan = samahanRange.getRange(\"After\");\nafterSamahan.load(\"start\");\nawait context.sync();\n\nconst periodResults = body.search(\".\", { matchCase: true });\nperiodResults.load(\"items,start\");\nawait context.sync();\n\nfor (const rng of periodResults.items) {\n  if (rng.start === afterSamahan.start) {\n    rng.font.bold = true;\n    break;\n  }\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Bold the possessive \"Ateneo's\" (curly apostrophe, matches the whole existing run)\n$r = $d.Content\n$r.Find.Execute(\"Ateneo\" + [char]0x2019 + \"s\") | Out-Null\n$r.Font.Bold = $true\n\n# 2) Bold \" ACMS (Automated Crowd Management System)\" but leave \" for IT Week\" unbolded\n$r = $d.Content\n$r.Find.Execute(\" ACMS (Automated Crowd Management System)\") | Out-Null\n$r.Font.Bold = $true\n\n# 3) Bold \"SAMAHAN System Development org\" (leave \" under \" before it unbolded)\n$r = $d.Content\n$r.Find.Execute(\"SAMAHAN System Development org\") | Out-Null\n$r.Font.Bold = $true\n$afterOrg = $r.Duplicate\n$afterOrg.Collapse(0)\n\n# 4) Bold the \".\" right after \"org\" but leave the trailing space unbolded\n$afterOrg.MoveEnd(1, 1) | Out-Null\nif ($afterOrg.Text -eq \".\") {\n    $afterOrg.Font.Bold = $true\n}\n\nWrite-Output \"done\"\n"}
